# Experiment #17 results (learn_rate 0.02 sweep) - replace the placeholder
# row 2 values with the real trial-0 results and append trials 1-4 (rows 3-6).
$arr = New-Object "object[,]" 5,20
$arr[0,0] = 0
$arr[0,1] = 0.0001
$arr[0,2] = 3
$arr[0,3] = 1000
$arr[0,4] = 200
$arr[0,5] = 4
$arr[0,6] = "('tanh', 'relu')"
$arr[0,7] = 100
$arr[0,8] = 0.02
$arr[0,9] = 1
$arr[0,10] = 0.9652000069618225
$arr[0,11] = 0.9287999868392944
$arr[0,12] = 77.307
$arr[0,13] = 0.0125
$arr[0,14] = 0.012
$arr[0,15] = 11
$arr[0,16] = 7.0279
$arr[0,17] = 0.9617999792098999
$arr[0,18] = 0.958899974822998
$arr[0,19] = 0.9623000025749207
$arr[1,0] = 1
$arr[1,1] = 0.0001
$arr[1,2] = 3
$arr[1,3] = 1000
$arr[1,4] = 200
$arr[1,5] = 4
$arr[1,6] = "('tanh', 'relu')"
$arr[1,7] = 100
$arr[1,8] = 0.02
$arr[1,9] = 2
$arr[1,10] = 0.9617999792098999
$arr[1,11] = 0.9182999730110168
$arr[1,12] = 43.371
$arr[1,13] = 0.0222
$arr[1,14] = 0.0212
$arr[1,15] = 6
$arr[1,16] = 7.2285
$arr[1,17] = 0.9567999839782715
$arr[1,18] = 0.9524999856948853
$arr[1,19] = 0.9546999931335449
$arr[2,0] = 2
$arr[2,1] = 0.0001
$arr[2,2] = 3
$arr[2,3] = 1000
$arr[2,4] = 200
$arr[2,5] = 4
$arr[2,6] = "('tanh', 'relu')"
$arr[2,7] = 100
$arr[2,8] = 0.02
$arr[2,9] = 3
$arr[2,10] = 0.9596999883651733
$arr[2,11] = 0.914900004863739
$arr[2,12] = 44.608
$arr[2,13] = 0.0215
$arr[2,14] = 0.0205
$arr[2,15] = 6
$arr[2,16] = 7.4347
$arr[2,17] = 0.9514999985694885
$arr[2,18] = 0.9526000022888184
$arr[2,19] = 0.9532999992370605
$arr[3,0] = 3
$arr[3,1] = 0.0001
$arr[3,2] = 3
$arr[3,3] = 1000
$arr[3,4] = 200
$arr[3,5] = 4
$arr[3,6] = "('tanh', 'relu')"
$arr[3,7] = 100
$arr[3,8] = 0.02
$arr[3,9] = 4
$arr[3,10] = 0.9646999835968018
$arr[3,11] = 0.9246000051498413
$arr[3,12] = 59.611
$arr[3,13] = 0.0162
$arr[3,14] = 0.0155
$arr[3,15] = 8
$arr[3,16] = 7.4514
$arr[3,17] = 0.9596999883651733
$arr[3,18] = 0.9584000110626221
$arr[3,19] = 0.9584000110626221
$arr[4,0] = 4
$arr[4,1] = 0.0001
$arr[4,2] = 3
$arr[4,3] = 1000
$arr[4,4] = 200
$arr[4,5] = 4
$arr[4,6] = "('tanh', 'relu')"
$arr[4,7] = 100
$arr[4,8] = 0.02
$arr[4,9] = 5
$arr[4,10] = 0.9589999914169312
$arr[4,11] = 0.9156000018119812
$arr[4,12] = 51.33
$arr[4,13] = 0.0187
$arr[4,14] = 0.0178
$arr[4,15] = 7
$arr[4,16] = 7.3329
$arr[4,17] = 0.9570000171661377
$arr[4,18] = 0.9545999765396118
$arr[4,19] = 0.954800009727478

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:T6")
$rng.Value = $arr

# Replicate the bordered/bold/centered style of A2 down to the newly added A3:A6 cells
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Edit complete"
